# 21 march Turkey data
# - Append row 12 (21 March 2020) to Sheet1
# - Add Sheet2 with tested/positive/percent_pos summary for 20 March 2020

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- Sheet1: new row 12 ----
# Copy formatting down from row 11 first so the new row inherits the same
# per-column styles (date format in C, bold in D) without creating new
# style entries.
$ws1.Range("A11:F11").Copy($ws1.Range("A12:F12"))

$ws1.Range("A12").Value = "Turkiye"
$ws1.Range("B12").Value = 11
$ws1.Range("C12").Value = 43911
$ws1.Range("D12").Value = 947
$ws1.Range("E12").Value = 21
$ws1.Range("F12").Formula = "=D12-E12"

$ws1.Range("F13").Select()

# ---- Sheet2: testing summary ----
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "Sheet2"

$ws2.Range("C1").Value = "positive"
$ws2.Range("B1").Value = "tested"
$ws2.Range("D1").Value = "percent_pos"

$ws1.Range("C11").Copy($ws2.Range("A2"))
$ws2.Range("A2").Value = 43910
$ws2.Range("B2").Value = 3656
$ws2.Range("C2").Value = 311
$ws2.Range("D2").Formula = "=C2/B2*100"

$ws2.Columns("A").ColumnWidth = 14.8

$ws2.Range("D3").Select()

$ws1.Activate()
